# Week5Religion.docx edit: update the `summary(iris.mis)` table values and
# the corresponding R console output (number of NAs / R-squares / nearest-
# neighbor imputation candidates table) to match a re-run of the script
# with a different random seed.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) The `kable(summary(iris.mis))` table (first table in the document).
#    Address cells directly (Row, Column) so each write is unambiguous
#    regardless of what the neighbouring cells currently contain.
# ---------------------------------------------------------------------
$t1 = $d.Tables(1)

$t1.Cell(2, 2).Range.Text = "Min. :4.300"
$t1.Cell(2, 3).Range.Text = "Min. :2.000"
$t1.Cell(2, 4).Range.Text = "Min. :1.000"
$t1.Cell(2, 6).Range.Text = "setosa :46"

$t1.Cell(3, 2).Range.Text = "1st Qu.:5.100"
$t1.Cell(3, 6).Range.Text = "versicolor:44"

$t1.Cell(4, 4).Range.Text = "Median :4.400"
$t1.Cell(4, 6).Range.Text = "virginica :40"

$t1.Cell(5, 2).Range.Text = "Mean :5.853"
$t1.Cell(5, 3).Range.Text = "Mean :3.055"
$t1.Cell(5, 4).Range.Text = "Mean :3.781"
$t1.Cell(5, 5).Range.Text = "Mean :1.227"
$t1.Cell(5, 6).Range.Text = "NA’s :20"

$t1.Cell(6, 3).Range.Text = "3rd Qu.:3.300"

$t1.Cell(7, 2).Range.Text = "Max. :7.700"

$t1.Cell(8, 2).Range.Text = "NA’s :13"
$t1.Cell(8, 3).Range.Text = "NA’s :13"
$t1.Cell(8, 4).Range.Text = "NA’s :15"
$t1.Cell(8, 5).Range.Text = "NA’s :14"

# ---------------------------------------------------------------------
# 2) "## Number of NAs:" console line - one number per column.
# ---------------------------------------------------------------------
$d.Content.Find.Execute( `
    "##           20           15           19           11           10 ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "##           13           13           15           14           20 ", 2)

# ---------------------------------------------------------------------
# 3) "R-squares for Predicting Non-Missing Values" console line.
# ---------------------------------------------------------------------
$d.Content.Find.Execute( `
    "##        0.880        0.677        0.977        0.952        0.988", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "##        0.817        0.676        0.983        0.963        0.988", 2)

# ---------------------------------------------------------------------
# 4) The `check` matrix dump (20 rows of nearest-neighbour candidates)
#    shrinks to 13 rows with new row labels / values. Update the first
#    13 lines in place, then delete the trailing 7 lines (each line is
#    its own run, preceded by a standalone run containing the line
#    break, so deleting a line also absorbs the break before it).
# ---------------------------------------------------------------------
$matrixUpdates = @(
    @{ Old = "## 6    5.4  5.5  6.0  5.4  5.7"; New = "## 15   5.1  5.0  5.1  5.1  5.4" },
    @{ Old = "## 9    4.6  5.0  4.9  4.8  5.0"; New = "## 19   5.4  5.0  4.8  4.9  5.5" },
    @{ Old = "## 14   4.4  4.8  4.8  5.0  5.0"; New = "## 31   5.1  4.7  4.8  5.1  5.1" },
    @{ Old = "## 22   5.4  5.8  5.0  5.4  5.3"; New = "## 38   5.0  4.7  5.0  5.1  5.2" },
    @{ Old = "## 28   5.2  5.0  4.8  5.1  5.0"; New = "## 50   4.8  4.8  4.8  5.5  4.9" },
    @{ Old = "## 38   5.1  4.6  5.4  5.0  5.4"; New = "## 60   6.7  6.7  6.4  6.9  6.3" },
    @{ Old = "## 39   4.6  4.6  5.0  4.8  5.0"; New = "## 61   5.4  4.8  5.1  5.4  4.8" },
    @{ Old = "## 42   4.9  4.9  4.8  4.8  4.9"; New = "## 80   5.1  5.5  5.5  5.2  5.4" },
    @{ Old = "## 44   4.8  5.1  5.4  5.1  4.6"; New = "## 91   5.6  5.7  6.2  5.7  5.9" },
    @{ Old = "## 50   5.1  5.2  5.1  4.6  4.6"; New = "## 104  6.3  6.5  6.4  6.7  6.3" },
    @{ Old = "## 51   6.0  5.7  6.4  6.5  6.0"; New = "## 105  6.7  5.8  6.0  6.2  6.4" },
    @{ Old = "## 52   6.7  5.9  6.7  6.7  6.1"; New = "## 132  7.7  7.7  7.7  7.7  7.7" },
    @{ Old = "## 55   6.7  6.7  6.1  5.6  5.8"; New = "## 147  6.7  6.3  6.0  6.2  5.7" }
)

foreach ($u in $matrixUpdates) {
    $d.Content.Find.Execute($u.Old, $true, $false, $false, $false, $false, $true, 1, $false, $u.New, 2)
}

$linesToRemove = @(
    "## 58   5.5  5.4  5.1  5.7  5.1",
    "## 71   6.9  6.1  6.1  5.9  6.1",
    "## 129  6.9  6.3  5.7  6.4  6.0",
    "## 134  6.9  6.3  6.1  6.4  6.7",
    "## 138  6.0  6.3  6.7  6.1  6.3",
    "## 141  6.7  6.0  6.4  6.0  6.3",
    "## 143  6.0  6.9  6.1  6.1  6.7"
)

foreach ($line in $linesToRemove) {
    $rng = $d.Content
    $rng.Find.Execute($line)
    # Pull the range start back one character to also grab the line-break
    # run immediately preceding this line, so the gap closes cleanly.
    $rng.Start = $rng.Start - 1
    $rng.Delete()
}

Write-Output "edit complete"
